# Replace imposter ticker data (rows 2-18, columns D-I) with GEHC data
# and let unused shared strings (the 16 extra tickers) be garbage collected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 62.71757793257134
$ws.Range("E2").Value = 58.11828994750977
$ws.Range("F2").Value = 65.98287283827436
$ws.Range("G2").Value = 52.76240683216319
$ws.Range("H2").Value = 456562037
$ws.Range("I2").Value = "GEHC"

$ws.Range("D3").Value = 62.71757793257134
$ws.Range("E3").Value = 58.11828994750977
$ws.Range("F3").Value = 65.98287283827436
$ws.Range("G3").Value = 52.76240683216319
$ws.Range("H3").Value = 456562037
$ws.Range("I3").Value = "GEHC"

$ws.Range("D4").Value = 62.71757793257134
$ws.Range("E4").Value = 58.11828994750977
$ws.Range("F4").Value = 65.98287283827436
$ws.Range("G4").Value = 52.76240683216319
$ws.Range("H4").Value = 456562037
$ws.Range("I4").Value = "GEHC"

$ws.Range("D5").Value = 62.71757793257134
$ws.Range("E5").Value = 58.11828994750977
$ws.Range("F5").Value = 65.98287283827436
$ws.Range("G5").Value = 52.76240683216319
$ws.Range("H5").Value = 456562037
$ws.Range("I5").Value = "GEHC"

$ws.Range("D6").Value = 62.71757793257134
$ws.Range("E6").Value = 58.11828994750977
$ws.Range("F6").Value = 65.98287283827436
$ws.Range("G6").Value = 52.76240683216319
$ws.Range("H6").Value = 456562037
$ws.Range("I6").Value = "GEHC"

$ws.Range("D7").Value = 62.71757793257134
$ws.Range("E7").Value = 58.11828994750977
$ws.Range("F7").Value = 65.98287283827436
$ws.Range("G7").Value = 52.76240683216319
$ws.Range("H7").Value = 456562037
$ws.Range("I7").Value = "GEHC"

$ws.Range("D8").Value = 60.40796957551981
$ws.Range("E8").Value = 69.20833587646484
$ws.Range("F8").Value = 73.61847632321393
$ws.Range("G8").Value = 53.26015748247527
$ws.Range("H8").Value = 456562037
$ws.Range("I8").Value = "GEHC"

$ws.Range("D9").Value = 81.63240967693162
$ws.Range("E9").Value = 80.97536468505859
$ws.Range("F9").Value = 87.43627672436963
$ws.Range("G9").Value = 77.21231091762122
$ws.Range("H9").Value = 456562037
$ws.Range("I9").Value = "GEHC"

$ws.Range("D10").Value = 81.07569570455831
$ws.Range("E10").Value = 77.67969512939453
$ws.Range("F10").Value = 82.80356398420423
$ws.Range("G10").Value = 75.93688145661964
$ws.Range("H10").Value = 456562037
$ws.Range("I10").Value = "GEHC"

$ws.Range("D11").Value = 67.74562725210113
$ws.Range("E11").Value = 66.32097625732422
$ws.Range("F11").Value = 68.73192175591446
$ws.Range("G11").Value = 62.11676111466002
$ws.Range("H11").Value = 456562037
$ws.Range("I11").Value = "GEHC"

$ws.Range("D12").Value = 76.3286334573708
$ws.Range("E12").Value = 73.11920166015625
$ws.Range("F12").Value = 77.40509025973188
$ws.Range("G12").Value = 71.07592759665707
$ws.Range("H12").Value = 456562037
$ws.Range("I12").Value = "GEHC"

$ws.Range("D13").Value = 91.04611682657868
$ws.Range("E13").Value = 76.01966857910156
$ws.Range("F13").Value = 91.19568485652604
$ws.Range("G13").Value = 75.86013492587813
$ws.Range("H13").Value = 456562037
$ws.Range("I13").Value = "GEHC"

$ws.Range("D14").Value = 78.15109021707227
$ws.Range("E14").Value = 84.41514587402344
$ws.Range("F14").Value = 88.02595846807378
$ws.Range("G14").Value = 75.48787041050669
$ws.Range("H14").Value = 456562037
$ws.Range("I14").Value = "GEHC"

$ws.Range("D15").Value = 93.61648922190386
$ws.Range("E15").Value = 87.1605224609375
$ws.Range("F15").Value = 93.61648922190386
$ws.Range("G15").Value = 84.71583994500412
$ws.Range("H15").Value = 456562037
$ws.Range("I15").Value = "GEHC"

$ws.Range("D16").Value = 78.7748299048564
$ws.Range("E16").Value = 88.13758087158203
$ws.Range("F16").Value = 89.63482171799377
$ws.Range("G16").Value = 77.57704027387132
$ws.Range("H16").Value = 456562037
$ws.Range("I16").Value = "GEHC"

$ws.Range("D17").Value = 80.25506632274828
$ws.Range("E17").Value = 70.22942352294922
$ws.Range("F17").Value = 80.72438778041541
$ws.Range("G17").Value = 57.5675567730737
$ws.Range("H17").Value = 456562037
$ws.Range("I17").Value = "GEHC"

$ws.Range("D18").Value = 73.75291612531082
$ws.Range("E18").Value = 71.25518798828125
$ws.Range("F18").Value = 78.69841478764484
$ws.Range("G18").Value = 70.25609673346942
$ws.Range("H18").Value = 456562037
$ws.Range("I18").Value = "GEHC"

